$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New rows for "lid" mutations (onderdeel/punt/lid cleanup)
# ---------------------------------------------------------------------
$ws.Range("F25").Value = "tweede tot en met achtste lid"
$ws.Range("F26").Value = "negende tot tiende lid"

# ---------------------------------------------------------------------
# 2. Text cleanup: "het eerste punt" -> "eerste punt"
# ---------------------------------------------------------------------
$ws.Range("E5").Value = "eerste punt"

# ---------------------------------------------------------------------
# 3. Colour-code the cells (fill formatting only, values untouched).
#    Re-use the workbook's existing "done" (green, same as A5) and
#    "pending" (orange, same as B6) highlight styles via PasteSpecial
#    so we don't create brand-new style/fill entries. (Multi-area Union
#    ranges only paste into their first area in this host, so each
#    contiguous block is pasted individually.)
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

# Green ("done") cells
$ws.Range("A5").Copy()
$ws.Range("E5").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Copy()
$ws.Range("E7").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Copy()
$ws.Range("F5:F17").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Copy()
$ws.Range("F23:F26").PasteSpecial($xlPasteFormats)

# Orange ("pending") cells
$ws.Range("B6").Copy()
$ws.Range("E6").PasteSpecial($xlPasteFormats)
$ws.Range("B6").Copy()
$ws.Range("E9").PasteSpecial($xlPasteFormats)
$ws.Range("B6").Copy()
$ws.Range("F18:F22").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Widen the new/visible "lid" (F) column.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 19.33

# ---------------------------------------------------------------------
# 5. Match the saved selection.
# ---------------------------------------------------------------------
$ws.Range("C23").Select()
